$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 1.53
$ws.Range("R2").Value = 2.5
$ws.Range("S2").Value = 1.9
$ws.Range("T2").Value = 2
$ws.Range("U2").Value = 2.38
$ws.Range("V2").Value = 1.57
$ws.Range("W2").Value = 1.25
$ws.Range("G3").Value = 1.95
$ws.Range("I3").Value = 4.2
$ws.Range("J3").Value = 2.75
$ws.Range("K3").Value = 1.95
$ws.Range("M3").Value = 1.1
$ws.Range("O3").Value = 1.5
$ws.Range("Q3").Value = 2.5
$ws.Range("R3").Value = 1.47
$ws.Range("V3").Value = 1.17
$ws.Range("AB3").Value = 8
$ws.Range("AD3").Value = 17
$ws.Range("AM3").Value = 19
$ws.Range("AO3").Value = 41
$ws.Range("AR3").Value = 1.87
$ws.Range("AS3").Value = 1.87
$ws.Range("G4").Value = 1.7
$ws.Range("I4").Value = 5.25
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 8
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 2.15
$ws.Range("R4").Value = 1.62
$ws.Range("U4").Value = 4.33
$ws.Range("AA4").Value = 5.5
$ws.Range("AC4").Value = 9
$ws.Range("AM4").Value = 26
$ws.Range("G5").Value = 4.75
$ws.Range("H5").Value = 3.75
$ws.Range("I5").Value = 1.67
$ws.Range("K5").Value = 2.3
$ws.Range("L5").Value = 2.25
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 13
$ws.Range("O5").Value = 1.22
$ws.Range("P5").Value = 4
$ws.Range("Q5").Value = 1.69
$ws.Range("R5").Value = 2.07
$ws.Range("S5").Value = 2
$ws.Range("T5").Value = 1.8
$ws.Range("U5").Value = 2.75
$ws.Range("V5").Value = 1.4
$ws.Range("W5").Value = 1.33
$ws.Range("X5").Value = 3.25
$ws.Range("Y5").Value = 1.73
$ws.Range("Z5").Value = 2
$ws.Range("AA5").Value = 15
$ws.Range("AB5").Value = 26
$ws.Range("AC5").Value = 15
$ws.Range("AF5").Value = 41
$ws.Range("AG5").Value = 13
$ws.Range("AH5").Value = 7.5
$ws.Range("AK5").Value = 201
$ws.Range("AL5").Value = 8
$ws.Range("AM5").Value = 8.5
$ws.Range("AP5").Value = 13
$ws.Range("AQ5").Value = 23
$ws.Range("G7").Value = 1.45
$ws.Range("H7").Value = 4.5
$ws.Range("I7").Value = 6
$ws.Range("L7").Value = 6
$ws.Range("M7").Value = 1.04
$ws.Range("N7").Value = 13
$ws.Range("Q7").Value = 1.7
$ws.Range("U7").Value = 2.63
$ws.Range("V7").Value = 1.44
$ws.Range("AB7").Value = 7.5
$ws.Range("AD7").Value = 11
$ws.Range("AF7").Value = 23
$ws.Range("AH7").Value = 9
$ws.Range("AK7").Value = 251
$ws.Range("AL7").Value = 15
$ws.Range("AM7").Value = 29
$ws.Range("AN7").Value = 17
